$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Logistic Regression)
$ws.Range("C2").Value = 0.6574201305393336
$ws.Range("D2").Value = 0.1924513915364087
$ws.Range("F2").Value = 0.296679400528945
$ws.Range("G2").Value = 0.7139790935411585
$ws.Range("H2").Value = 0.6530078465562337
$ws.Range("I2").Value = 0.1849315068493151
$ws.Range("J2").Value = 0.6639344262295082
$ws.Range("K2").Value = 0.2892857142857143
$ws.Range("L2").Value = 0.7156744639422338

# Row 3 (Decision Tree)
$ws.Range("C3").Value = 0.6264456658651093
$ws.Range("D3").Value = 0.1814636863728153
$ws.Range("F3").Value = 0.2854170089261267
$ws.Range("G3").Value = 0.6995266135535857
$ws.Range("H3").Value = 0.6232311716182684
$ws.Range("I3").Value = 0.172939649578196
$ws.Range("J3").Value = 0.6721311475409836
$ws.Range("K3").Value = 0.2750967741935484
$ws.Range("L3").Value = 0.6956156708914069

# Row 4 (Random Forest)
$ws.Range("C4").Value = 0.5980762624527654
$ws.Range("D4").Value = 0.1771934292627022
$ws.Range("E4").Value = 0.7135897435897436
$ws.Range("F4").Value = 0.2838926859124757
$ws.Range("G4").Value = 0.7095110630217939
$ws.Range("H4").Value = 0.5917108175172692
$ws.Range("I4").Value = 0.1689705882352941
$ws.Range("J4").Value = 0.7244640605296343
$ws.Range("K4").Value = 0.2740281421416647
$ws.Range("L4").Value = 0.7055389205264639
